$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data for columns B, C, D (Circulation, ILL Loans, ILL Borrows) for rows 3-59
$data = @(
    @(11206,1834,1972),
    @(5732,720,765),
    @(18905,1757,1800),
    @(226,200,35),
    @(12027,2090,1933),
    @(1522,352,371),
    @(1555,298,201),
    @(512,123,73),
    @(39,65,14),
    @(3,0,0),
    @(265,72,106),
    @(841,366,200),
    @(1498,555,237),
    @(741,352,96),
    @(501,249,45),
    @(4066,831,856),
    @(463,208,168),
    @(4937,575,899),
    @(79,107,6),
    @(4369,636,884),
    @(347,117,39),
    @(4496,817,995),
    @(20866,2075,2607),
    @(1445,515,299),
    @(0,0,0),
    @(1357,284,287),
    @(353,152,72),
    @(4116,700,706),
    @(141,63,73),
    @(635,512,78),
    @(4063,1055,732),
    @(2660,794,704),
    @(1646,178,390),
    @(14304,1632,1467),
    @(1946,895,327),
    @(6032,566,752),
    @(271,302,99),
    @(454,167,149),
    @(878,133,68),
    @(4265,250,135),
    @(154,33,52),
    @(332,17,30),
    @(594,14,1),
    @(780,263,134),
    @(2931,974,616),
    @(8233,1189,1383),
    @(3202,1004,377),
    @(3079,347,477),
    @(7961,864,1250),
    @(1116,152,308),
    @(3815,867,719),
    @(408,307,181),
    @(568,416,52),
    @(753,181,224),
    @(3345,1345,639),
    @(6223,447,286),
    @(173688,29225,26444)
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    $ws.Cells.Item($r, 2).Value = $rowVals[0]
    $ws.Cells.Item($r, 3).Value = $rowVals[1]
    $ws.Cells.Item($r, 4).Value = $rowVals[2]
}

# Restore the active cell selection on the frozen-pane sheet view to A2
$ws.Range("A2").Select()
